$d = $word.ActiveDocument

# Update the date line (unique text in the document, safe to use a document-wide Find/Replace)
$d.Content.Find.Execute("2025-02-12 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-02-13 Thursday", 2)

# Helper to set the text of a specific table cell while preserving its
# run/paragraph formatting (replaces only the visible text, not the
# trailing cell-mark character).
function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

$t = $d.Tables.Item(1)

# Row 1 (table rows 1,5,9,13,17 contain the visible data; the rows in
# between are blank spacer rows)
Set-CellText $t 1 1 "16÷8=2, 0"
Set-CellText $t 1 2 "68÷8=8, 4"
Set-CellText $t 1 3 "21÷7=3, 0"
Set-CellText $t 1 4 "34÷9=3, 7"
Set-CellText $t 1 5 "50÷7=7, 1"

# Row 5
Set-CellText $t 5 1 "84÷4=21, 0"
Set-CellText $t 5 2 "10÷2=5, 0"
Set-CellText $t 5 3 "74÷2=37, 0"
Set-CellText $t 5 4 "12÷6=2, 0"
Set-CellText $t 5 5 "38÷9=4, 2"

# Row 9
Set-CellText $t 9 1 "85÷6=14, 1"
Set-CellText $t 9 2 "45÷7=6, 3"
Set-CellText $t 9 3 "64÷8=8, 0"
Set-CellText $t 9 4 "25÷5=5, 0"
Set-CellText $t 9 5 "44÷3=14, 2"

# Row 13
Set-CellText $t 13 1 "86÷8=10, 6"
Set-CellText $t 13 2 "37÷9=4, 1"
Set-CellText $t 13 3 "29÷3=9, 2"
Set-CellText $t 13 4 "61÷2=30, 1"
Set-CellText $t 13 5 "56÷3=18, 2"

# Row 17
Set-CellText $t 17 1 "76÷9=8, 4"
Set-CellText $t 17 2 "65÷9=7, 2"
Set-CellText $t 17 3 "56÷3=18, 2"
Set-CellText $t 17 4 "44÷6=7, 2"
Set-CellText $t 17 5 "19÷5=3, 4"

Write-Host "All replacements complete"
